# Insert a new record row at row 110 (shifting the existing rows 110-189
# down to 111-190) and populate it with the new price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(110).Insert()

$ws.Cells.Item(110, 1).Value  = 10
$ws.Cells.Item(110, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(110, 3).Value  = "La Araucanía"
$ws.Cells.Item(110, 4).Value  = "2023-07-19"
$ws.Cells.Item(110, 5).Value  = 9
$ws.Cells.Item(110, 6).Value  = 100112035
$ws.Cells.Item(110, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(110, 8).Value  = "Sin especificar"
$ws.Cells.Item(110, 9).Value  = "Primera"
$ws.Cells.Item(110, 10).Value = 50
$ws.Cells.Item(110, 11).Value = 25000
$ws.Cells.Item(110, 12).Value = 25000
$ws.Cells.Item(110, 13).Value = 25000
$ws.Cells.Item(110, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(110, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(110, 16).Value = 1667
$ws.Cells.Item(110, 17).Value = 15
$ws.Cells.Item(110, 18).Value = "Hortaliza"
